# [PHOENIX-5914] changes in Council Management
# Update the "citizen" test-data row (row 26) on the registeredUserDetails
# sheet with a new mobile number / password, and move the view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: dataName="citizen", id=9036544535 -> 9916534408, password=akhi2506 -> kurnool_eGov@123
$ws.Cells.Item(26, 2).Value = "9916534408"

$ws.Cells.Item(26, 3).Value = "kurnool_eGov@123"
$ws.Cells.Item(26, 3).NumberFormat = "@"

# Move the scroll position / active selection on the sheet view.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("C28").Select()
